$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-tracker sheet: append the newest scrape as a new row right after
# the current last row (row 37 -> new row 38), keeping columns A-D
# (Date, Price, Discount, Incredible) as plain text, matching every
# existing row in the sheet.
$newRow = 38

# Leading "'" forces each value to be stored as literal text instead of
# being auto-coerced (dates -> serial numbers, digit strings -> numbers).
$ws.Cells.Item($newRow, 1).Value = "'2026-02-07"
$ws.Cells.Item($newRow, 2).Value = "'43990000"
$ws.Cells.Item($newRow, 3).Value = "'0"
$ws.Cells.Item($newRow, 4).Value = "'0"

# Re-normalize the style so the new cells use the same default/general
# style as every other row (drops the "quote prefix" look only, keeps
# the text content).
$ws.Range("A$newRow`:D$newRow").Style = "Normal"
